$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.654881412372163
$ws.Range("B3").Value = 1.587511678723217
$ws.Range("B4").Value = 1.547709735660874
$ws.Range("B5").Value = 1.531880906409128
$ws.Range("B6").Value = 1.529276077001441
$ws.Range("B7").Value = 1.547494683280945
$ws.Range("B8").Value = 1.631326407921563
$ws.Range("B9").Value = 1.808243947042229
$ws.Range("B10").Value = 1.94604602191572
$ws.Range("B11").Value = 2.01047653750777
$ws.Range("B12").Value = 2.035128669714027
$ws.Range("B13").Value = 2.029808060759422
$ws.Range("B14").Value = 2.012499581582972
$ws.Range("B15").Value = 2.001930766140447
$ws.Range("B16").Value = 1.941870673897654
$ws.Range("B17").Value = 1.905474422068664
$ws.Range("B18").Value = 1.884704199872317
$ws.Range("B19").Value = 1.877699827562765
$ws.Range("B20").Value = 1.909331876987835
$ws.Range("B21").Value = 2.017576592263879
$ws.Range("B22").Value = 2.089800933539834
$ws.Range("B23").Value = 2.05111700721784
$ws.Range("B24").Value = 1.907587440708198
$ws.Range("B25").Value = 1.759023393854648

$ws.Range("C2").Value = 0.3269421160261743
$ws.Range("C3").Value = 0.3113641504571092
$ws.Range("C4").Value = 0.3021723861943997
$ws.Range("C5").Value = 0.2985197732467952
$ws.Range("C6").Value = 0.2979188619854369
$ws.Range("C7").Value = 0.3021227497480652
$ws.Range("C8").Value = 0.3214929308547312
$ws.Range("C9").Value = 0.3624750792876057
$ws.Range("C10").Value = 0.3944671454913475
$ws.Range("C11").Value = 0.4094424976639175
$ws.Range("C12").Value = 0.415174936307352
$ws.Range("C13").Value = 0.4139375986545986
$ws.Range("C14").Value = 0.4099128689963436
$ws.Range("C15").Value = 0.407455655316852
$ws.Range("C16").Value = 0.393497040321904
$ws.Range("C17").Value = 0.3850426090880887
$ws.Range("C18").Value = 0.380219508272404
$ws.Range("C19").Value = 0.3785932741720046
$ws.Range("C20").Value = 0.3859384860678858
$ws.Range("C21").Value = 0.4110933504914556
$ws.Range("C22").Value = 0.4278929910618103
$ws.Range("C23").Value = 0.4188935049162978
$ws.Range("C24").Value = 0.3855333434610202
$ws.Range("C25").Value = 0.3510618189815489

$ws.Range("D2").Value = 0.04116341890146913
$ws.Range("D3").Value = 0.03808527330614453
$ws.Range("D4").Value = 0.03620435605166961
$ws.Range("D5").Value = 0.03543994555655416
$ws.Range("D6").Value = 0.03531313537883563
$ws.Range("D7").Value = 0.0361940387971984
$ws.Range("D8").Value = 0.04010006669525268
$ws.Range("D9").Value = 0.0478413228032224
$ws.Range("D10").Value = 0.05359221828011584
$ws.Range("D11").Value = 0.05622521823757154
$ws.Range("D12").Value = 0.05722493247141358
$ws.Range("D13").Value = 0.05700950445243791
$ws.Range("D14").Value = 0.05630741083118096
$ws.Range("D15").Value = 0.05587771062643299
$ws.Range("D16").Value = 0.05342050462174086
$ws.Range("D17").Value = 0.05191757373357575
$ws.Range("D18").Value = 0.05105470363831444
$ws.Range("D19").Value = 0.05076281449564135
$ws.Range("D20").Value = 0.05207739855512727
$ws.Range("D21").Value = 0.05651355889820309
$ws.Range("D22").Value = 0.05942842401695714
$ws.Range("D23").Value = 0.05787120500767173
$ws.Range("D24").Value = 0.0520051380799913
$ws.Range("D25").Value = 0.04573697145535505

$ws.Range("E2").Value = 0.07090668655138455
$ws.Range("E3").Value = 0.0715131709158543
$ws.Range("E4").Value = 0.07191249689968959
$ws.Range("E5").Value = 0.07208201892930965
$ws.Range("E6").Value = 0.07211057879045768
$ws.Range("E7").Value = 0.07191475560147431
$ws.Range("E8").Value = 0.07111022167740177
$ws.Range("E9").Value = 0.06974548337809239
$ws.Range("E10").Value = 0.06887151593473551
$ws.Range("E11").Value = 0.06850163978818546
$ws.Range("E12").Value = 0.06836554187147126
$ws.Range("E13").Value = 0.06839467685753053
$ws.Range("E14").Value = 0.06849036353677906
$ws.Range("E15").Value = 0.06854949041495484
$ws.Range("E16").Value = 0.06889624410985462
$ws.Range("E17").Value = 0.06911604807206917
$ws.Range("E18").Value = 0.06924508135163787
$ws.Range("E19").Value = 0.06928921821067657
$ws.Range("E20").Value = 0.06909237978234728
$ws.Range("E21").Value = 0.06846215052964055
$ws.Range("E22").Value = 0.06807336955784127
$ws.Range("E23").Value = 0.06827876007995215
$ws.Range("E24").Value = 0.06910307191381904
$ws.Range("E25").Value = 0.07009199948940825

$ws.Range("F2").Value = 5.035716064856729
$ws.Range("F3").Value = 4.851349761971932
$ws.Range("F4").Value = 4.739545318433045
$ws.Range("F5").Value = 4.694327541323219
$ws.Range("F6").Value = 4.686839677559135
$ws.Range("F7").Value = 4.738934115475132
$ws.Range("F8").Value = 4.971851739380128
$ws.Range("F9").Value = 5.440070794699466
$ws.Range("F10").Value = 5.791629124467505
$ws.Range("F11").Value = 5.953329208850903
$ws.Range("F12").Value = 6.014825675715713
$ws.Range("F13").Value = 6.001569418181248
$ws.Range("F14").Value = 5.958383214121568
$ws.Range("F15").Value = 5.931965087199785
$ws.Range("F16").Value = 5.781098115288358
$ws.Range("F17").Value = 5.689006943400869
$ws.Range("F18").Value = 5.636205009110938
$ws.Range("F19").Value = 5.618355501492204
$ws.Range("F20").Value = 5.698792911956588
$ws.Range("F21").Value = 5.97106080764371
$ws.Range("F22").Value = 6.150547587691278
$ws.Range("F23").Value = 6.054607832565921
$ws.Range("F24").Value = 5.694368232372682
$ws.Range("F25").Value = 5.312119290078272

$ws.Range("J2").Value = 0.1825092295514565
$ws.Range("J3").Value = 0.1800861087913717
$ws.Range("J4").Value = 0.1786455481354849
$ws.Range("J5").Value = 0.1780702718662894
$ws.Range("J6").Value = 0.1779754548219969
$ws.Range("J7").Value = 0.1786377422794558
$ws.Range("J8").Value = 0.1816638604472587
$ws.Range("J9").Value = 0.1879786103686101
$ws.Range("J10").Value = 0.1928585155890943
$ws.Range("J11").Value = 0.1951326180401409
$ws.Range("J12").Value = 0.1960017025120351
$ws.Range("J13").Value = 0.1958141746856086
$ws.Range("J14").Value = 0.1952039584005263
$ws.Range("J15").Value = 0.1948312202676732
$ws.Range("J16").Value = 0.1927110000422729
$ws.Range("J17").Value = 0.1914242938552917
$ws.Range("J18").Value = 0.1906893086201933
$ws.Range("J19").Value = 0.1904413260672513
$ws.Range("J20").Value = 0.1915607377251192
$ws.Range("J21").Value = 0.1953829773677072
$ws.Range("J22").Value = 0.1979273261512091
$ws.Range("J23").Value = 0.1965650777904102
$ws.Range("J24").Value = 0.1914990366377083
$ws.Range("J25").Value = 0.1862286961929982

$ws.Range("M2").Value = 0.4906818933942247
$ws.Range("M3").Value = 0.4803755351054306
$ws.Range("M4").Value = 0.4745060722008105
$ws.Range("M5").Value = 0.4722291538982617
$ws.Range("M6").Value = 0.4718580036487552
$ws.Range("M7").Value = 0.4744749000345223
$ws.Range("M8").Value = 0.4870328087104596
$ws.Range("M9").Value = 0.5153194574814037
$ws.Range("M10").Value = 0.5383669275127403
$ws.Range("M11").Value = 0.5493515098151676
$ws.Range("M12").Value = 0.5535836053970584
$ws.Range("M13").Value = 0.5526689156966214
$ws.Range("M14").Value = 0.5496982311915062
$ws.Range("M15").Value = 0.5478880558615415
$ws.Range("M16").Value = 0.5376591703948606
$ws.Range("M17").Value = 0.5315125704641375
$ws.Range("M18").Value = 0.5280242533333137
$ws.Range("M19").Value = 0.5268512326574424
$ws.Range("M20").Value = 0.5321620136780183
$ws.Range("M21").Value = 0.5505688218807663
$ws.Range("M22").Value = 0.5630214304542918
$ws.Range("M23").Value = 0.5563363756703623
$ws.Range("M24").Value = 0.5318682588910377
$ws.Range("M25").Value = 0.5072718644532017
